# Eng Test Data - 12 Mar 2025
# Update the AddOpportunity sheet's "Sector" sample value and refresh the
# workbook's active-sheet / selection state to match the authored edit.

$wb = $excel.ActiveWorkbook

# --- AddOpportunity sheet: replace the AdditionalClient/Sector sample text ---
$ws = $wb.Worksheets.Item("AddOpportunity")

$cell = $ws.Range("E2")
$cell.Value = "CSDN-0000001546"

# New data is longer than the old "Dealership & Rental Services" label, so
# give it wrap text + vertically centered alignment (matches the new style
# added to the workbook).
$cell.WrapText = $true
$cell.VerticalAlignment = -4108

# --- Users sheet keeps its own selection, but is no longer the active tab ---
$users = $wb.Worksheets.Item("Users")
$users.Range("E15").Select()

# --- AddOpportunity becomes the active sheet with a fresh selection ---
$ws.Activate()
$ws.Range("E5").Select()
